$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$blankStyle = $ws.Range("D4").Style

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $blankStyle
}

Set-TextValue 'D2' '62.978.62'
Set-TextValue 'E2' '  +2.78%  '
Set-TextValue 'D3' '2.950.53'
Set-TextValue 'E3' '  +0.94%  '
Set-TextValue 'D5' '594.91'
Set-TextValue 'E5' '  -0.29%  '
Set-TextValue 'D6' '148.30'
Set-TextValue 'E6' '  +2.22%  '
Set-TextValue 'D7' '1.00'
Set-TextValue 'E7' '  -0.02%  '
Set-TextValue 'D8' '2.948.79'
Set-TextValue 'E8' '  +0.96%  '
Set-TextValue 'E9' '  +1.26%  '
Set-TextValue 'D10' '7.28'
Set-TextValue 'E10' '  +4.44%  '
Set-TextValue 'E11' '  +6.87%  '
Set-TextValue 'D12' '0.442'
Set-TextValue 'E12' '  +1.05%  '
Set-TextValue 'E13' '  +5.53%  '
Set-TextValue 'D14' '32.80'
Set-TextValue 'E14' '  -1.79%  '
Set-TextValue 'E15' '  -0.60%  '
Set-TextValue 'D16' '3.435.85'
Set-TextValue 'E16' '  +0.79%  '
Set-TextValue 'D17' '62.880.53'
Set-TextValue 'E17' '  +2.61%  '
Set-TextValue 'D18' '6.72'
Set-TextValue 'E18' '  +0.55%  '
Set-TextValue 'D19' '2.945.67'
Set-TextValue 'E19' '  +0.74%  '
Set-TextValue 'D20' '442.58'
Set-TextValue 'E20' '  +2.74%  '
Set-TextValue 'D21' '13.43'
Set-TextValue 'E21' '  -0.20%  '
Set-TextValue 'D22' '0.668'
Set-TextValue 'E22' '  -1.10%  '
Set-TextValue 'D23' '7.05'
Set-TextValue 'E23' '  -0.05%  '
Set-TextValue 'D24' '11.29'
Set-TextValue 'E24' '  +4.41%  '
Set-TextValue 'D25' '81.05'
Set-TextValue 'E25' '  -0.91%  '
Set-TextValue 'D26' '2.15'
Set-TextValue 'E26' '  -1.35%  '
Set-TextValue 'D27' '11.82'
Set-TextValue 'E27' '  +0.98%  '
Set-TextValue 'E28' '  -0.02%  '
Set-TextValue 'D29' '2.21'
Set-TextValue 'E29' '  +0.34%  '
Set-TextValue 'D30' '7.20'
Set-TextValue 'E30' '  +4.54%  '
Set-TextValue 'D31' '2.62'
Set-TextValue 'E31' '  +0.35%  '
Set-TextValue 'D32' '0.0000103'
Set-TextValue 'E32' '  +17.22%  '
Set-TextValue 'B33' 'Hedera'
Set-TextValue 'C33' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D33' '0.109'
Set-TextValue 'E33' '  -0.46%  '
Set-TextValue 'B34' 'EthereumClassic'
Set-TextValue 'C34' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D34' '26.43'
Set-TextValue 'E34' '  -0.39%  '
Set-TextValue 'D35' '1.00'
Set-TextValue 'E35' '  -0.06%  '
Set-TextValue 'E36' '  -1.67%  '
Set-TextValue 'D37' '3.15'
Set-TextValue 'E37' '  +5.93%  '
Set-TextValue 'D38' '5.60'
Set-TextValue 'E38' '  -0.04%  '
Set-TextValue 'D39' '49.69'
Set-TextValue 'E39' '  -0.39%  '
Set-TextValue 'D40' '2.04'
Set-TextValue 'E40' '  +2.52%  '
Set-TextValue 'D41' '8.51'
Set-TextValue 'E41' '  -0.47%  '
Set-TextValue 'D42' '0.118'
Set-TextValue 'E42' '  -3.74%  '
Set-TextValue 'D43' '0.281'
Set-TextValue 'E43' '  +0.05%  '
Set-TextValue 'D44' '38.99'
Set-TextValue 'E44' '  -7.61%  '
Set-TextValue 'B45' 'Monero'
Set-TextValue 'C45' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D45' '135.46'
Set-TextValue 'E45' '  +1.27%  '
Set-TextValue 'B46' 'Maker'
Set-TextValue 'C46' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D46' '2.694.01'
Set-TextValue 'E46' '  -0.12%  '
Set-TextValue 'D47' '0.0337'
Set-TextValue 'E47' '  -2.27%  '
Set-TextValue 'D48' '362.62'
Set-TextValue 'E48' '  +0.57%  '
Set-TextValue 'E49' '  -0.02%  '
Set-TextValue 'E50' '  -0.40%  '
Set-TextValue 'D51' '22.78'
Set-TextValue 'E51' '  -3.04%  '
